# Atualizando cronograma: mover o marcador "X" de progresso para a direita
# nas linhas 6 e 7 (colunas J/K -> L, e J -> M respectivamente).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Linha 6: estava em J6, passa para L6 (K6 permanece vazio)
$ws.Range("J6").Value = ""
$ws.Range("K6").Value = ""
$ws.Range("L6").Value = "X"

# Linha 7: estava em J7, passa para M7 (K7 e L7 permanecem vazios)
$ws.Range("J7").Value = ""
$ws.Range("M7").Value = "X"

# Atualiza a celula selecionada/ativa para refletir a ultima edicao
$ws.Range("M7").Select()
